$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: propagate existing cell *formatting* (style indexes) to the cells
# that will hold the new layout, using Copy + PasteSpecial(xlPasteFormats)
# so the existing style entries in styles.xml are reused instead of new
# (duplicate) styles being minted. The order below is dependency-safe: every
# source cell still holds the required style at the time it is read.
# ---------------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$ws.Range("E1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("D1").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("D1").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("D1").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: clear the soon-to-be-deleted column E contents (its data either
# moved elsewhere or is dropped) before we write the new grid values. Also
# fully clear D6:D11 (value + style) since those "CONSIDERACIONES" rows no
# longer exist past row 5 in the new layout - a plain ClearContents would
# leave a dangling, styled-but-empty cell behind.
# ---------------------------------------------------------------------------
$ws.Range("E1:E12").ClearContents()
$ws.Range("D6:D11").Clear()

# ---------------------------------------------------------------------------
# Step 3: write the new header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "camino_optico"
$ws.Range("B1").Value = "θ_grados"
$ws.Range("C1").Value = " θ_grados_minutos"
$ws.Range("D1").Value = "CONSIDERACIONES"

# ---------------------------------------------------------------------------
# Step 4: write the new data grid (A: camino optico values, B: angle in
# degrees, C: angle degrees/minutes label, D: considerations column).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 8.9861003778057107
$ws.Range("B2").Value = 103.58329999999999
$ws.Range("C2").Value = "103°35'"
$ws.Range("D2").Value = "500ml agua"

$ws.Range("A3").Value = 8.8741196746494193
$ws.Range("B3").Value = 104.667
$ws.Range("C3").Value = "104°40'"
$ws.Range("D3").Value = "Concentracion=64.2209gms"

$ws.Range("A4").Value = 8.7749643873921208
$ws.Range("B4").Value = 105.25
$ws.Range("C4").Value = "105°15'"
$ws.Range("D4").Value = "lamda = 532nm"

$ws.Range("A5").Value = 8.6162636914152007
$ws.Range("B5").Value = 105.5
$ws.Range("C5").Value = "105°30'"
$ws.Range("D5").Value = 7.3

$ws.Range("A6").Value = 8.4118963379252296
$ws.Range("B6").Value = 105.833
$ws.Range("C6").Value = "105°50'"

$ws.Range("A7").Value = 8.2042671823874596
$ws.Range("B7").Value = 106.417
$ws.Range("C7").Value = "106°25'"

$ws.Range("A8").Value = 7.9063265806567804
$ws.Range("B8").Value = 106.667
$ws.Range("C8").Value = "106°40'"

$ws.Range("A9").Value = 7.48331477354788
$ws.Range("B9").Value = 107.5
$ws.Range("C9").Value = "107°30'"

$ws.Range("A10").Value = 6.8818602136341003
$ws.Range("B10").Value = 108.667
$ws.Range("C10").Value = "108°40'"

$ws.Range("A11").Value = 6.11882341631134
$ws.Range("B11").Value = 109.5
$ws.Range("C11").Value = "109°30'"

# ---------------------------------------------------------------------------
# Step 5: drop the now-empty column E entirely and tidy up the sheet view /
# column widths to match the new A:D layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).Delete()

$ws.Columns.Item(1).ColumnWidth = 14.88671875

$ws.Range("C16").Select()
